$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric columns (H): comma-decimal/thousands -> dot-decimal text ---
# Force text format so Excel keeps these as literal strings (not coerced
# to numbers, which would drop the formatting/trailing zeros), then reset
# the cell style back to Normal so no stray number format lingers on it.
$c = $ws.Range('H2')
$c.NumberFormat = '@'
$c.Value = '720.00'
$c.Style = 'Normal'
$c = $ws.Range('H3')
$c.NumberFormat = '@'
$c.Value = '3020.00'
$c.Style = 'Normal'
$c = $ws.Range('H4')
$c.NumberFormat = '@'
$c.Value = '12438.00'
$c.Style = 'Normal'
$c = $ws.Range('H5')
$c.NumberFormat = '@'
$c.Value = '158.00'
$c.Style = 'Normal'
$c = $ws.Range('H6')
$c.NumberFormat = '@'
$c.Value = '23697.53'
$c.Style = 'Normal'
$c = $ws.Range('H7')
$c.NumberFormat = '@'
$c.Value = '919.98'
$c.Style = 'Normal'
$c = $ws.Range('H8')
$c.NumberFormat = '@'
$c.Value = '3090.00'
$c.Style = 'Normal'
$c = $ws.Range('H9')
$c.NumberFormat = '@'
$c.Value = '37139.22'
$c.Style = 'Normal'
$c = $ws.Range('H10')
$c.NumberFormat = '@'
$c.Value = '66558.50'
$c.Style = 'Normal'
$c = $ws.Range('H11')
$c.NumberFormat = '@'
$c.Value = '5750.00'
$c.Style = 'Normal'
$c = $ws.Range('H12')
$c.NumberFormat = '@'
$c.Value = '2881.20'
$c.Style = 'Normal'
$c = $ws.Range('H13')
$c.NumberFormat = '@'
$c.Value = '42470.02'
$c.Style = 'Normal'
$c = $ws.Range('H14')
$c.NumberFormat = '@'
$c.Value = '14723.20'
$c.Style = 'Normal'
$c = $ws.Range('H15')
$c.NumberFormat = '@'
$c.Value = '745.98'
$c.Style = 'Normal'
$c = $ws.Range('H16')
$c.NumberFormat = '@'
$c.Value = '24869.95'
$c.Style = 'Normal'
$c = $ws.Range('H17')
$c.NumberFormat = '@'
$c.Value = '10528.76'
$c.Style = 'Normal'
$c = $ws.Range('H18')
$c.NumberFormat = '@'
$c.Value = '117.00'
$c.Style = 'Normal'
$c = $ws.Range('H19')
$c.NumberFormat = '@'
$c.Value = '3060.00'
$c.Style = 'Normal'
$c = $ws.Range('H20')
$c.NumberFormat = '@'
$c.Value = '541.64'
$c.Style = 'Normal'
$c = $ws.Range('H21')
$c.NumberFormat = '@'
$c.Value = '3377.70'
$c.Style = 'Normal'
$c = $ws.Range('H22')
$c.NumberFormat = '@'
$c.Value = '60.72'
$c.Style = 'Normal'
$c = $ws.Range('H23')
$c.NumberFormat = '@'
$c.Value = '175.00'
$c.Style = 'Normal'
$c = $ws.Range('H24')
$c.NumberFormat = '@'
$c.Value = '900.00'
$c.Style = 'Normal'
$c = $ws.Range('H25')
$c.NumberFormat = '@'
$c.Value = '39556.00'
$c.Style = 'Normal'
$c = $ws.Range('H26')
$c.NumberFormat = '@'
$c.Value = '7931.39'
$c.Style = 'Normal'
$c = $ws.Range('H27')
$c.NumberFormat = '@'
$c.Value = '190.00'
$c.Style = 'Normal'
$c = $ws.Range('H28')
$c.NumberFormat = '@'
$c.Value = '835.30'
$c.Style = 'Normal'
$c = $ws.Range('H29')
$c.NumberFormat = '@'
$c.Value = '38.34'
$c.Style = 'Normal'
$c = $ws.Range('H30')
$c.NumberFormat = '@'
$c.Value = '110.00'
$c.Style = 'Normal'
$c = $ws.Range('H31')
$c.NumberFormat = '@'
$c.Value = '8865.06'
$c.Style = 'Normal'
$c = $ws.Range('H32')
$c.NumberFormat = '@'
$c.Value = '199.76'
$c.Style = 'Normal'
$c = $ws.Range('H33')
$c.NumberFormat = '@'
$c.Value = '14.41'
$c.Style = 'Normal'
$c = $ws.Range('H34')
$c.NumberFormat = '@'
$c.Value = '217.00'
$c.Style = 'Normal'
$c = $ws.Range('H35')
$c.NumberFormat = '@'
$c.Value = '13104.50'
$c.Style = 'Normal'
$c = $ws.Range('H36')
$c.NumberFormat = '@'
$c.Value = '3130.00'
$c.Style = 'Normal'
$c = $ws.Range('H37')
$c.NumberFormat = '@'
$c.Value = '15600.00'
$c.Style = 'Normal'
$c = $ws.Range('H38')
$c.NumberFormat = '@'
$c.Value = '12309.00'
$c.Style = 'Normal'
$c = $ws.Range('H39')
$c.NumberFormat = '@'
$c.Value = '240.00'
$c.Style = 'Normal'
$c = $ws.Range('H40')
$c.NumberFormat = '@'
$c.Value = '7486.57'
$c.Style = 'Normal'
$c = $ws.Range('H41')
$c.NumberFormat = '@'
$c.Value = '20.00'
$c.Style = 'Normal'
$c = $ws.Range('H42')
$c.NumberFormat = '@'
$c.Value = '465.68'
$c.Style = 'Normal'
$c = $ws.Range('H43')
$c.NumberFormat = '@'
$c.Value = '410.00'
$c.Style = 'Normal'
$c = $ws.Range('H44')
$c.NumberFormat = '@'
$c.Value = '180.16'
$c.Style = 'Normal'
$c = $ws.Range('H45')
$c.NumberFormat = '@'
$c.Value = '2489.00'
$c.Style = 'Normal'
$c = $ws.Range('H46')
$c.NumberFormat = '@'
$c.Value = '613.30'
$c.Style = 'Normal'
$c = $ws.Range('H47')
$c.NumberFormat = '@'
$c.Value = '34.72'
$c.Style = 'Normal'
$c = $ws.Range('H48')
$c.NumberFormat = '@'
$c.Value = '152.00'
$c.Style = 'Normal'
$c = $ws.Range('H49')
$c.NumberFormat = '@'
$c.Value = '61.73'
$c.Style = 'Normal'
$c = $ws.Range('H50')
$c.NumberFormat = '@'
$c.Value = '73.90'
$c.Style = 'Normal'
$c = $ws.Range('H51')
$c.NumberFormat = '@'
$c.Value = '1540.15'
$c.Style = 'Normal'
$c = $ws.Range('H52')
$c.NumberFormat = '@'
$c.Value = '3616.84'
$c.Style = 'Normal'
$c = $ws.Range('H53')
$c.NumberFormat = '@'
$c.Value = '1063.30'
$c.Style = 'Normal'
$c = $ws.Range('H54')
$c.NumberFormat = '@'
$c.Value = '2053.00'
$c.Style = 'Normal'
$c = $ws.Range('H55')
$c.NumberFormat = '@'
$c.Value = '548.00'
$c.Style = 'Normal'
$c = $ws.Range('H56')
$c.NumberFormat = '@'
$c.Value = '1216.43'
$c.Style = 'Normal'
$c = $ws.Range('H57')
$c.NumberFormat = '@'
$c.Value = '1520.00'
$c.Style = 'Normal'
$c = $ws.Range('H58')
$c.NumberFormat = '@'
$c.Value = '4630.20'
$c.Style = 'Normal'
$c = $ws.Range('H59')
$c.NumberFormat = '@'
$c.Value = '160.00'
$c.Style = 'Normal'
$c = $ws.Range('H60')
$c.NumberFormat = '@'
$c.Value = '1466.74'
$c.Style = 'Normal'
$c = $ws.Range('H61')
$c.NumberFormat = '@'
$c.Value = '286.50'
$c.Style = 'Normal'
$c = $ws.Range('H62')
$c.NumberFormat = '@'
$c.Value = '396.04'
$c.Style = 'Normal'
$c = $ws.Range('H63')
$c.NumberFormat = '@'
$c.Value = '10072.68'
$c.Style = 'Normal'
$c = $ws.Range('H64')
$c.NumberFormat = '@'
$c.Value = '150.00'
$c.Style = 'Normal'
$c = $ws.Range('H65')
$c.NumberFormat = '@'
$c.Value = '7200.00'
$c.Style = 'Normal'
$c = $ws.Range('H66')
$c.NumberFormat = '@'
$c.Value = '3500.00'
$c.Style = 'Normal'
$c = $ws.Range('H67')
$c.NumberFormat = '@'
$c.Value = '5200.00'
$c.Style = 'Normal'
$c = $ws.Range('H68')
$c.NumberFormat = '@'
$c.Value = '4352.37'
$c.Style = 'Normal'
$c = $ws.Range('H69')
$c.NumberFormat = '@'
$c.Value = '252.00'
$c.Style = 'Normal'
$c = $ws.Range('H70')
$c.NumberFormat = '@'
$c.Value = '1000.00'
$c.Style = 'Normal'
$c = $ws.Range('H71')
$c.NumberFormat = '@'
$c.Value = '1902.00'
$c.Style = 'Normal'
$c = $ws.Range('H72')
$c.NumberFormat = '@'
$c.Value = '854.00'
$c.Style = 'Normal'
$c = $ws.Range('H73')
$c.NumberFormat = '@'
$c.Value = '925.00'
$c.Style = 'Normal'
$c = $ws.Range('H74')
$c.NumberFormat = '@'
$c.Value = '415.00'
$c.Style = 'Normal'
$c = $ws.Range('H75')
$c.NumberFormat = '@'
$c.Value = '2340.88'
$c.Style = 'Normal'
$c = $ws.Range('H76')
$c.NumberFormat = '@'
$c.Value = '2572.00'
$c.Style = 'Normal'
$c = $ws.Range('H77')
$c.NumberFormat = '@'
$c.Value = '700.00'
$c.Style = 'Normal'
$c = $ws.Range('H78')
$c.NumberFormat = '@'
$c.Value = '86.20'
$c.Style = 'Normal'
$c = $ws.Range('H79')
$c.NumberFormat = '@'
$c.Value = '260000.00'
$c.Style = 'Normal'
$c = $ws.Range('H80')
$c.NumberFormat = '@'
$c.Value = '61589.50'
$c.Style = 'Normal'
$c = $ws.Range('H81')
$c.NumberFormat = '@'
$c.Value = '300.00'
$c.Style = 'Normal'
$c = $ws.Range('H82')
$c.NumberFormat = '@'
$c.Value = '25.81'
$c.Style = 'Normal'
$c = $ws.Range('H83')
$c.NumberFormat = '@'
$c.Value = '406.29'
$c.Style = 'Normal'
$c = $ws.Range('H84')
$c.NumberFormat = '@'
$c.Value = '136.20'
$c.Style = 'Normal'
$c = $ws.Range('H85')
$c.NumberFormat = '@'
$c.Value = '5968.60'
$c.Style = 'Normal'
$c = $ws.Range('H86')
$c.NumberFormat = '@'
$c.Value = '150.00'
$c.Style = 'Normal'
$c = $ws.Range('H87')
$c.NumberFormat = '@'
$c.Value = '178.00'
$c.Style = 'Normal'
$c = $ws.Range('H88')
$c.NumberFormat = '@'
$c.Value = '7274.28'
$c.Style = 'Normal'
$c = $ws.Range('H89')
$c.NumberFormat = '@'
$c.Value = '105.00'
$c.Style = 'Normal'
$c = $ws.Range('H90')
$c.NumberFormat = '@'
$c.Value = '252.00'
$c.Style = 'Normal'
$c = $ws.Range('H91')
$c.NumberFormat = '@'
$c.Value = '165.00'
$c.Style = 'Normal'
$c = $ws.Range('H92')
$c.NumberFormat = '@'
$c.Value = '83.70'
$c.Style = 'Normal'
$c = $ws.Range('H93')
$c.NumberFormat = '@'
$c.Value = '6604.25'
$c.Style = 'Normal'
$c = $ws.Range('H94')
$c.NumberFormat = '@'
$c.Value = '1863.86'
$c.Style = 'Normal'
$c = $ws.Range('H95')
$c.NumberFormat = '@'
$c.Value = '25.02'
$c.Style = 'Normal'
$c = $ws.Range('H96')
$c.NumberFormat = '@'
$c.Value = '1182.95'
$c.Style = 'Normal'
$c = $ws.Range('H97')
$c.NumberFormat = '@'
$c.Value = '14.99'
$c.Style = 'Normal'
$c = $ws.Range('H98')
$c.NumberFormat = '@'
$c.Value = '348.00'
$c.Style = 'Normal'
$c = $ws.Range('H99')
$c.NumberFormat = '@'
$c.Value = '2152.00'
$c.Style = 'Normal'
$c = $ws.Range('H100')
$c.NumberFormat = '@'
$c.Value = '237.51'
$c.Style = 'Normal'
$c = $ws.Range('H101')
$c.NumberFormat = '@'
$c.Value = '60.00'
$c.Style = 'Normal'
$c = $ws.Range('H102')
$c.NumberFormat = '@'
$c.Value = '2000.00'
$c.Style = 'Normal'
$c = $ws.Range('H103')
$c.NumberFormat = '@'
$c.Value = '810.00'
$c.Style = 'Normal'
$c = $ws.Range('H104')
$c.NumberFormat = '@'
$c.Value = '2310.00'
$c.Style = 'Normal'
$c = $ws.Range('H105')
$c.NumberFormat = '@'
$c.Value = '920.00'
$c.Style = 'Normal'
$c = $ws.Range('H106')
$c.NumberFormat = '@'
$c.Value = '30876.00'
$c.Style = 'Normal'
$c = $ws.Range('H107')
$c.NumberFormat = '@'
$c.Value = '25680.00'
$c.Style = 'Normal'
$c = $ws.Range('H108')
$c.NumberFormat = '@'
$c.Value = '180.00'
$c.Style = 'Normal'
$c = $ws.Range('H109')
$c.NumberFormat = '@'
$c.Value = '550.00'
$c.Style = 'Normal'
$c = $ws.Range('H110')
$c.NumberFormat = '@'
$c.Value = '1200.00'
$c.Style = 'Normal'
$c = $ws.Range('H111')
$c.NumberFormat = '@'
$c.Value = '1680.30'
$c.Style = 'Normal'
$c = $ws.Range('H112')
$c.NumberFormat = '@'
$c.Value = '4892.46'
$c.Style = 'Normal'
$c = $ws.Range('H113')
$c.NumberFormat = '@'
$c.Value = '1260.00'
$c.Style = 'Normal'
$c = $ws.Range('H114')
$c.NumberFormat = '@'
$c.Value = '1492.50'
$c.Style = 'Normal'
$c = $ws.Range('H115')
$c.NumberFormat = '@'
$c.Value = '111.60'
$c.Style = 'Normal'
$c = $ws.Range('H116')
$c.NumberFormat = '@'
$c.Value = '387.00'
$c.Style = 'Normal'
$c = $ws.Range('H117')
$c.NumberFormat = '@'
$c.Value = '3245.10'
$c.Style = 'Normal'
$c = $ws.Range('H118')
$c.NumberFormat = '@'
$c.Value = '98250.00'
$c.Style = 'Normal'
$c = $ws.Range('H119')
$c.NumberFormat = '@'
$c.Value = '122640.00'
$c.Style = 'Normal'
$c = $ws.Range('H120')
$c.NumberFormat = '@'
$c.Value = '17868.61'
$c.Style = 'Normal'
$c = $ws.Range('H121')
$c.NumberFormat = '@'
$c.Value = '1000.00'
$c.Style = 'Normal'
$c = $ws.Range('H122')
$c.NumberFormat = '@'
$c.Value = '700.00'
$c.Style = 'Normal'
$c = $ws.Range('H123')
$c.NumberFormat = '@'
$c.Value = '2400.00'
$c.Style = 'Normal'
$c = $ws.Range('H124')
$c.NumberFormat = '@'
$c.Value = '250.00'
$c.Style = 'Normal'
$c = $ws.Range('H125')
$c.NumberFormat = '@'
$c.Value = '2210.00'
$c.Style = 'Normal'
$c = $ws.Range('H126')
$c.NumberFormat = '@'
$c.Value = '384.00'
$c.Style = 'Normal'
$c = $ws.Range('H127')
$c.NumberFormat = '@'
$c.Value = '600.00'
$c.Style = 'Normal'
$c = $ws.Range('H128')
$c.NumberFormat = '@'
$c.Value = '350.00'
$c.Style = 'Normal'
$c = $ws.Range('H129')
$c.NumberFormat = '@'
$c.Value = '2500.00'
$c.Style = 'Normal'
$c = $ws.Range('H130')
$c.NumberFormat = '@'
$c.Value = '8000.00'
$c.Style = 'Normal'
$c = $ws.Range('H131')
$c.NumberFormat = '@'
$c.Value = '23364.90'
$c.Style = 'Normal'
$c = $ws.Range('H132')
$c.NumberFormat = '@'
$c.Value = '1500.00'
$c.Style = 'Normal'
$c = $ws.Range('H133')
$c.NumberFormat = '@'
$c.Value = '950.00'
$c.Style = 'Normal'
$c = $ws.Range('H134')
$c.NumberFormat = '@'
$c.Value = '750.00'
$c.Style = 'Normal'
$c = $ws.Range('H135')
$c.NumberFormat = '@'
$c.Value = '7400.00'
$c.Style = 'Normal'
$c = $ws.Range('H136')
$c.NumberFormat = '@'
$c.Value = '3334.27'
$c.Style = 'Normal'
$c = $ws.Range('H137')
$c.NumberFormat = '@'
$c.Value = '1500.00'
$c.Style = 'Normal'
$c = $ws.Range('H138')
$c.NumberFormat = '@'
$c.Value = '200.00'
$c.Style = 'Normal'
$c = $ws.Range('H139')
$c.NumberFormat = '@'
$c.Value = '660.00'
$c.Style = 'Normal'
$c = $ws.Range('H140')
$c.NumberFormat = '@'
$c.Value = '8480.00'
$c.Style = 'Normal'
$c = $ws.Range('H141')
$c.NumberFormat = '@'
$c.Value = '1000.00'
$c.Style = 'Normal'
$c = $ws.Range('H142')
$c.NumberFormat = '@'
$c.Value = '270.00'
$c.Style = 'Normal'
$c = $ws.Range('H143')
$c.NumberFormat = '@'
$c.Value = '350.00'
$c.Style = 'Normal'
$c = $ws.Range('H144')
$c.NumberFormat = '@'
$c.Value = '125.00'
$c.Style = 'Normal'
$c = $ws.Range('H145')
$c.NumberFormat = '@'
$c.Value = '1800.00'
$c.Style = 'Normal'
$c = $ws.Range('H146')
$c.NumberFormat = '@'
$c.Value = '12500.00'
$c.Style = 'Normal'
$c = $ws.Range('H147')
$c.NumberFormat = '@'
$c.Value = '350.00'
$c.Style = 'Normal'
$c = $ws.Range('H148')
$c.NumberFormat = '@'
$c.Value = '6570.00'
$c.Style = 'Normal'
$c = $ws.Range('H149')
$c.NumberFormat = '@'
$c.Value = '33202.50'
$c.Style = 'Normal'
$c = $ws.Range('H150')
$c.NumberFormat = '@'
$c.Value = '356.00'
$c.Style = 'Normal'
$c = $ws.Range('H151')
$c.NumberFormat = '@'
$c.Value = '2096.00'
$c.Style = 'Normal'
$c = $ws.Range('H152')
$c.NumberFormat = '@'
$c.Value = '2855.00'
$c.Style = 'Normal'
$c = $ws.Range('H153')
$c.NumberFormat = '@'
$c.Value = '870.00'
$c.Style = 'Normal'
$c = $ws.Range('H154')
$c.NumberFormat = '@'
$c.Value = '1578.34'
$c.Style = 'Normal'
$c = $ws.Range('H155')
$c.NumberFormat = '@'
$c.Value = '10811.00'
$c.Style = 'Normal'
$c = $ws.Range('H156')
$c.NumberFormat = '@'
$c.Value = '2925.00'
$c.Style = 'Normal'
$c = $ws.Range('H157')
$c.NumberFormat = '@'
$c.Value = '205.00'
$c.Style = 'Normal'
$c = $ws.Range('H158')
$c.NumberFormat = '@'
$c.Value = '325.00'
$c.Style = 'Normal'
$c = $ws.Range('H159')
$c.NumberFormat = '@'
$c.Value = '1700.00'
$c.Style = 'Normal'
$c = $ws.Range('H160')
$c.NumberFormat = '@'
$c.Value = '110.00'
$c.Style = 'Normal'
$c = $ws.Range('H161')
$c.NumberFormat = '@'
$c.Value = '10931.40'
$c.Style = 'Normal'
$c = $ws.Range('H162')
$c.NumberFormat = '@'
$c.Value = '269.68'
$c.Style = 'Normal'
$c = $ws.Range('H163')
$c.NumberFormat = '@'
$c.Value = '155.00'
$c.Style = 'Normal'
$c = $ws.Range('H164')
$c.NumberFormat = '@'
$c.Value = '1100.00'
$c.Style = 'Normal'
$c = $ws.Range('H165')
$c.NumberFormat = '@'
$c.Value = '250.00'
$c.Style = 'Normal'
$c = $ws.Range('H166')
$c.NumberFormat = '@'
$c.Value = '298.56'
$c.Style = 'Normal'
$c = $ws.Range('H167')
$c.NumberFormat = '@'
$c.Value = '3375.63'
$c.Style = 'Normal'
$c = $ws.Range('H168')
$c.NumberFormat = '@'
$c.Value = '4262.66'
$c.Style = 'Normal'
$c = $ws.Range('H169')
$c.NumberFormat = '@'
$c.Value = '544.49'
$c.Style = 'Normal'
$c = $ws.Range('H170')
$c.NumberFormat = '@'
$c.Value = '664.00'
$c.Style = 'Normal'
$c = $ws.Range('H171')
$c.NumberFormat = '@'
$c.Value = '2977.00'
$c.Style = 'Normal'
$c = $ws.Range('H172')
$c.NumberFormat = '@'
$c.Value = '590.00'
$c.Style = 'Normal'
$c = $ws.Range('H173')
$c.NumberFormat = '@'
$c.Value = '2625.00'
$c.Style = 'Normal'
$c = $ws.Range('H174')
$c.NumberFormat = '@'
$c.Value = '6570.00'
$c.Style = 'Normal'
$c = $ws.Range('H175')
$c.NumberFormat = '@'
$c.Value = '1147.55'
$c.Style = 'Normal'
$c = $ws.Range('H176')
$c.NumberFormat = '@'
$c.Value = '148.93'
$c.Style = 'Normal'
$c = $ws.Range('H177')
$c.NumberFormat = '@'
$c.Value = '1560.00'
$c.Style = 'Normal'
$c = $ws.Range('H178')
$c.NumberFormat = '@'
$c.Value = '3822.97'
$c.Style = 'Normal'
$c = $ws.Range('H179')
$c.NumberFormat = '@'
$c.Value = '1893.20'
$c.Style = 'Normal'
$c = $ws.Range('H180')
$c.NumberFormat = '@'
$c.Value = '8728.25'
$c.Style = 'Normal'
$c = $ws.Range('H181')
$c.NumberFormat = '@'
$c.Value = '179.90'
$c.Style = 'Normal'
$c = $ws.Range('H182')
$c.NumberFormat = '@'
$c.Value = '4408.00'
$c.Style = 'Normal'
$c = $ws.Range('H183')
$c.NumberFormat = '@'
$c.Value = '6136.80'
$c.Style = 'Normal'
$c = $ws.Range('H184')
$c.NumberFormat = '@'
$c.Value = '946.00'
$c.Style = 'Normal'
$c = $ws.Range('H185')
$c.NumberFormat = '@'
$c.Value = '774.00'
$c.Style = 'Normal'
$c = $ws.Range('H186')
$c.NumberFormat = '@'
$c.Value = '146287.42'
$c.Style = 'Normal'
$c = $ws.Range('H187')
$c.NumberFormat = '@'
$c.Value = '3640.00'
$c.Style = 'Normal'
$c = $ws.Range('H188')
$c.NumberFormat = '@'
$c.Value = '70000.00'
$c.Style = 'Normal'
$c = $ws.Range('H189')
$c.NumberFormat = '@'
$c.Value = '1200.00'
$c.Style = 'Normal'
$c = $ws.Range('H190')
$c.NumberFormat = '@'
$c.Value = '1197387.98'
$c.Style = 'Normal'
$c = $ws.Range('H191')
$c.NumberFormat = '@'
$c.Value = '1200.00'
$c.Style = 'Normal'
$c = $ws.Range('H192')
$c.NumberFormat = '@'
$c.Value = '1080.00'
$c.Style = 'Normal'
$c = $ws.Range('H193')
$c.NumberFormat = '@'
$c.Value = '33800.00'
$c.Style = 'Normal'
$c = $ws.Range('H194')
$c.NumberFormat = '@'
$c.Value = '7410.28'
$c.Style = 'Normal'
$c = $ws.Range('H195')
$c.NumberFormat = '@'
$c.Value = '20650.00'
$c.Style = 'Normal'
$c = $ws.Range('H196')
$c.NumberFormat = '@'
$c.Value = '136500.00'
$c.Style = 'Normal'
$c = $ws.Range('H197')
$c.NumberFormat = '@'
$c.Value = '62500.00'
$c.Style = 'Normal'
$c = $ws.Range('H198')
$c.NumberFormat = '@'
$c.Value = '37000.00'
$c.Style = 'Normal'
$c = $ws.Range('H199')
$c.NumberFormat = '@'
$c.Value = '20000.00'
$c.Style = 'Normal'
$c = $ws.Range('H200')
$c.NumberFormat = '@'
$c.Value = '35000.00'
$c.Style = 'Normal'
$c = $ws.Range('H201')
$c.NumberFormat = '@'
$c.Value = '234201.00'
$c.Style = 'Normal'
$c = $ws.Range('H202')
$c.NumberFormat = '@'
$c.Value = '32000.00'
$c.Style = 'Normal'
$c = $ws.Range('H203')
$c.NumberFormat = '@'
$c.Value = '111500.00'
$c.Style = 'Normal'
$c = $ws.Range('H204')
$c.NumberFormat = '@'
$c.Value = '211799.00'
$c.Style = 'Normal'
$c = $ws.Range('H205')
$c.NumberFormat = '@'
$c.Value = '130000.00'
$c.Style = 'Normal'
$c = $ws.Range('H206')
$c.NumberFormat = '@'
$c.Value = '391005.45'
$c.Style = 'Normal'
$c = $ws.Range('H207')
$c.NumberFormat = '@'
$c.Value = '153334.75'
$c.Style = 'Normal'
$c = $ws.Range('H208')
$c.NumberFormat = '@'
$c.Value = '756.25'
$c.Style = 'Normal'
$c = $ws.Range('H209')
$c.NumberFormat = '@'
$c.Value = '66247.50'
$c.Style = 'Normal'
$c = $ws.Range('H210')
$c.NumberFormat = '@'
$c.Value = '27225.00'
$c.Style = 'Normal'
$c = $ws.Range('H211')
$c.NumberFormat = '@'
$c.Value = '9500.00'
$c.Style = 'Normal'
$c = $ws.Range('H212')
$c.NumberFormat = '@'
$c.Value = '26236.00'
$c.Style = 'Normal'
$c = $ws.Range('H213')
$c.NumberFormat = '@'
$c.Value = '2800.00'
$c.Style = 'Normal'
$c = $ws.Range('H214')
$c.NumberFormat = '@'
$c.Value = '1150.00'
$c.Style = 'Normal'
$c = $ws.Range('H215')
$c.NumberFormat = '@'
$c.Value = '4747.00'
$c.Style = 'Normal'
$c = $ws.Range('H216')
$c.NumberFormat = '@'
$c.Value = '1545.00'
$c.Style = 'Normal'

# --- Proper-name columns (E/F): comma separators -> period separators ---
$ws.Range('E29').Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range('F29').Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range('E73').Value = 'FERNANDEZ. MARIO HUGO'
$ws.Range('E75').Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range('F75').Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range('E76').Value = 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'
$ws.Range('E89').Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range('F89').Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range('E91').Value = 'TRABICHET MARIA. VERGARA ADEL Y OTRA'
$ws.Range('F91').Value = 'TRABICHET MARIA. VERGARA ADEL Y OTRA'
$ws.Range('E102').Value = 'RICCOTTI. MARIANA EDITH'
$ws.Range('E142').Value = 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'
$ws.Range('E145').Value = 'DODERA. JORGE ABELARDO'
$ws.Range('E151').Value = 'ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN'
$ws.Range('E159').Value = 'DODERA. JORGE ABELARDO'
$ws.Range('E163').Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range('F163').Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range('E172').Value = 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'
$ws.Range('E189').Value = 'DODERA. JORGE ABELARDO'
